# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G (K) values for rows 2-16 with freshly recalculated counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 3
    10 = 1
    11 = 1
    12 = 3
    13 = 2
    14 = 1
    15 = 1
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
